$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume change (E) cells with latest crypto data refresh.
# Each value is textual (not numeric) in the source data, so we temporarily force
# a text number format before assigning, then restore the cells original style
# afterwards so no visible formatting change is introduced.

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.912.35'
$ws.Range('D2').Style = $origStyle

$origStyle = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('E2').Style = $origStyle

$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.815.50'
$ws.Range('D3').Style = $origStyle

$origStyle = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E3').Style = $origStyle

$origStyle = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('E4').Style = $origStyle

$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.07'
$ws.Range('D5').Style = $origStyle

$origStyle = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E5').Style = $origStyle

$origStyle = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E6').Style = $origStyle

$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4642'
$ws.Range('D7').Style = $origStyle

$origStyle = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('E7').Style = $origStyle

$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3662'
$ws.Range('D8').Style = $origStyle

$origStyle = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.26%  '
$ws.Range('E8').Style = $origStyle

$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07355'
$ws.Range('D9').Style = $origStyle

$origStyle = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('E9').Style = $origStyle

$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8712'
$ws.Range('D10').Style = $origStyle

$origStyle = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('E10').Style = $origStyle

$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.27'
$ws.Range('D11').Style = $origStyle

$origStyle = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E11').Style = $origStyle

$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.824.09'
$ws.Range('D12').Style = $origStyle

$origStyle = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('E12').Style = $origStyle

$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.370'
$ws.Range('D13').Style = $origStyle

$origStyle = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('E13').Style = $origStyle

$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07098'
$ws.Range('D14').Style = $origStyle

$origStyle = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('E14').Style = $origStyle

$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.506'
$ws.Range('D15').Style = $origStyle

$origStyle = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('E15').Style = $origStyle

$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.45'
$ws.Range('D16').Style = $origStyle

$origStyle = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('E16').Style = $origStyle

$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('D17').Style = $origStyle

$origStyle = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('E17').Style = $origStyle

$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008718'
$ws.Range('D18').Style = $origStyle

$origStyle = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('E18').Style = $origStyle

$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = $origStyle

$origStyle = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('E19').Style = $origStyle

$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.67'
$ws.Range('D20').Style = $origStyle

$origStyle = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E20').Style = $origStyle

$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.934.27'
$ws.Range('D21').Style = $origStyle

$origStyle = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E21').Style = $origStyle

$origStyle = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E22').Style = $origStyle

$origStyle = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E23').Style = $origStyle

$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.057.52'
$ws.Range('D24').Style = $origStyle

$origStyle = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('E24').Style = $origStyle

$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.893'
$ws.Range('D25').Style = $origStyle

$origStyle = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('E25').Style = $origStyle

$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.82'
$ws.Range('D26').Style = $origStyle

$origStyle = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('E26').Style = $origStyle

$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.33'
$ws.Range('D27').Style = $origStyle

$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.133'
$ws.Range('D28').Style = $origStyle

$origStyle = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.75%  '
$ws.Range('E28').Style = $origStyle

$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.266'
$ws.Range('D29').Style = $origStyle

$origStyle = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E29').Style = $origStyle

$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.39'
$ws.Range('D30').Style = $origStyle

$origStyle = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('E30').Style = $origStyle

$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08910'
$ws.Range('D31').Style = $origStyle

$origStyle = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E31').Style = $origStyle

$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7557'
$ws.Range('D32').Style = $origStyle

$origStyle = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('E32').Style = $origStyle

$origStyle = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('E33').Style = $origStyle

$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.479'
$ws.Range('D34').Style = $origStyle

$origStyle = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('E34').Style = $origStyle

$origStyle = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('E35').Style = $origStyle

$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = $origStyle

$origStyle = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('E36').Style = $origStyle

$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.087'
$ws.Range('D37').Style = $origStyle

$origStyle = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('E37').Style = $origStyle

$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05282'
$ws.Range('D38').Style = $origStyle

$origStyle = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('E38').Style = $origStyle

$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.992'
$ws.Range('D39').Style = $origStyle

$origStyle = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.28%  '
$ws.Range('E39').Style = $origStyle

$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01950'
$ws.Range('D40').Style = $origStyle

$origStyle = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('E40').Style = $origStyle

$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.250'
$ws.Range('D41').Style = $origStyle

$origStyle = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('E41').Style = $origStyle

$origStyle = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E42').Style = $origStyle

$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.323'
$ws.Range('D43').Style = $origStyle

$origStyle = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('E43').Style = $origStyle

$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1652'
$ws.Range('D44').Style = $origStyle

$origStyle = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('E44').Style = $origStyle

$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.439'
$ws.Range('D45').Style = $origStyle

$origStyle = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('E45').Style = $origStyle

$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4871'
$ws.Range('D46').Style = $origStyle

$origStyle = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.92%  '
$ws.Range('E46').Style = $origStyle

$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.40'
$ws.Range('D47').Style = $origStyle

$origStyle = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('E47').Style = $origStyle

$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = $origStyle

$origStyle = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('E48').Style = $origStyle

$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.35'
$ws.Range('D49').Style = $origStyle

$origStyle = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('E49').Style = $origStyle

$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.661'
$ws.Range('D50').Style = $origStyle

$origStyle = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('E50').Style = $origStyle

$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06296'
$ws.Range('D51').Style = $origStyle

$origStyle = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.04%  '
$ws.Range('E51').Style = $origStyle

